$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.333947
$ws.Cells.Item(2, 8).Value = 1.001841
$ws.Cells.Item(2, 9).Value = 0.9184595666969813
$ws.Cells.Item(2, 10).Value = 0.9184595666969813
$ws.Cells.Item(2, 13).Value = 10.50827633333333
$ws.Cells.Item(2, 14).Value = 31.524829
$ws.Cells.Item(2, 15).Value = 0.1682660991018133
$ws.Cells.Item(2, 16).Value = 0.1682660991018134
$ws.Cells.Item(2, 17).Value = 3.509207356687666
$ws.Cells.Item(2, 18).Value = 31.582866210189
$ws.Cells.Item(2, 19).Value = 0.1545456084708428
$ws.Cells.Item(2, 20).Value = 0.1545456084708428
# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.333947
$ws.Cells.Item(3, 8).Value = 1.001841
$ws.Cells.Item(3, 9).Value = 0.9184595666969813
$ws.Cells.Item(3, 10).Value = 0.9184595666969813
$ws.Cells.Item(3, 15).Value = 0.4955285863849104
$ws.Cells.Item(3, 16).Value = 0.4955285863849105
$ws.Cells.Item(3, 17).Value = 10.33430126491967
$ws.Cells.Item(3, 18).Value = 93.00871138427699
$ws.Cells.Item(3, 19).Value = 0.4551229707370525
$ws.Cells.Item(3, 20).Value = 0.4551229707370525
# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.333947
$ws.Cells.Item(4, 8).Value = 1.001841
$ws.Cells.Item(4, 9).Value = 0.9184595666969813
$ws.Cells.Item(4, 10).Value = 0.9184595666969813
$ws.Cells.Item(4, 13).Value = 6.495209666666667
$ws.Cells.Item(4, 14).Value = 19.485629
$ws.Cells.Item(4, 15).Value = 0.1040059814559238
$ws.Cells.Item(4, 16).Value = 0.1040059814559238
$ws.Cells.Item(4, 17).Value = 2.169055782554334
$ws.Cells.Item(4, 18).Value = 19.521502042989
$ws.Cells.Item(4, 19).Value = 0.09552528866190202
$ws.Cells.Item(4, 20).Value = 0.09552528866190202
# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.333947
$ws.Cells.Item(5, 8).Value = 1.001841
$ws.Cells.Item(5, 9).Value = 0.9184595666969813
$ws.Cells.Item(5, 10).Value = 0.9184595666969813
$ws.Cells.Item(5, 13).Value = 9.909791666666667
$ws.Cells.Item(5, 14).Value = 29.729375
$ws.Cells.Item(5, 15).Value = 0.1586827309986352
$ws.Cells.Item(5, 16).Value = 0.1586827309986352
$ws.Cells.Item(5, 17).Value = 3.309345197708333
$ws.Cells.Item(5, 18).Value = 29.784106779375
$ws.Cells.Item(5, 19).Value = 0.1457436723553001
$ws.Cells.Item(5, 20).Value = 0.1457436723553001
# Row 6
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.333947
$ws.Cells.Item(6, 8).Value = 1.001841
$ws.Cells.Item(6, 9).Value = 0.9184595666969813
$ws.Cells.Item(6, 10).Value = 0.9184595666969813
$ws.Cells.Item(6, 13).Value = 4.591137333333333
$ws.Cells.Item(6, 14).Value = 13.773412
$ws.Cells.Item(6, 15).Value = 0.07351660205871713
$ws.Cells.Item(6, 16).Value = 0.07351660205871713
$ws.Cells.Item(6, 17).Value = 1.533196539054667
$ws.Cells.Item(6, 18).Value = 13.798768851492
$ws.Cells.Item(6, 19).Value = 0.06752202647188374
$ws.Cells.Item(6, 20).Value = 0.06752202647188374
# Row 7
$ws.Cells.Item(7, 9).Value = 0.08154043330301874
$ws.Cells.Item(7, 10).Value = 0.08154043330301874
$ws.Cells.Item(7, 13).Value = 10.50827633333333
$ws.Cells.Item(7, 14).Value = 31.524829
$ws.Cells.Item(7, 15).Value = 0.1682660991018133
$ws.Cells.Item(7, 16).Value = 0.1682660991018134
$ws.Cells.Item(7, 17).Value = 0.3115458739718889
$ws.Cells.Item(7, 18).Value = 2.803912865747
$ws.Cells.Item(7, 19).Value = 0.01372049063097055
$ws.Cells.Item(7, 20).Value = 0.01372049063097056
# Row 8
$ws.Cells.Item(8, 9).Value = 0.08154043330301874
$ws.Cells.Item(8, 10).Value = 0.08154043330301874
$ws.Cells.Item(8, 15).Value = 0.4955285863849104
$ws.Cells.Item(8, 16).Value = 0.4955285863849105
$ws.Cells.Item(8, 19).Value = 0.04040561564785795
$ws.Cells.Item(8, 20).Value = 0.04040561564785795
# Row 9
$ws.Cells.Item(9, 9).Value = 0.08154043330301874
$ws.Cells.Item(9, 10).Value = 0.08154043330301874
$ws.Cells.Item(9, 13).Value = 6.495209666666667
$ws.Cells.Item(9, 14).Value = 19.485629
$ws.Cells.Item(9, 15).Value = 0.1040059814559238
$ws.Cells.Item(9, 16).Value = 0.1040059814559238
$ws.Cells.Item(9, 17).Value = 0.1925678111274444
$ws.Cells.Item(9, 18).Value = 1.733110300147
$ws.Cells.Item(9, 19).Value = 0.008480692794021758
$ws.Cells.Item(9, 20).Value = 0.008480692794021758
# Row 10
$ws.Cells.Item(10, 9).Value = 0.08154043330301874
$ws.Cells.Item(10, 10).Value = 0.08154043330301874
$ws.Cells.Item(10, 13).Value = 9.909791666666667
$ws.Cells.Item(10, 14).Value = 29.729375
$ws.Cells.Item(10, 15).Value = 0.1586827309986352
$ws.Cells.Item(10, 16).Value = 0.1586827309986352
$ws.Cells.Item(10, 17).Value = 0.2938022000694445
$ws.Cells.Item(10, 18).Value = 2.644219800625
$ws.Cells.Item(10, 19).Value = 0.01293905864333508
$ws.Cells.Item(10, 20).Value = 0.01293905864333508
# Row 11
$ws.Cells.Item(11, 9).Value = 0.08154043330301874
$ws.Cells.Item(11, 10).Value = 0.08154043330301874
$ws.Cells.Item(11, 13).Value = 4.591137333333333
$ws.Cells.Item(11, 14).Value = 13.773412
$ws.Cells.Item(11, 15).Value = 0.07351660205871713
$ws.Cells.Item(11, 16).Value = 0.07351660205871713
$ws.Cells.Item(11, 17).Value = 0.1361165092795555
$ws.Cells.Item(11, 18).Value = 1.225048583516
$ws.Cells.Item(11, 19).Value = 0.005994575586833395
$ws.Cells.Item(11, 20).Value = 0.005994575586833395

Write-Output "Applied TPM updates"